# Insert two new weekly-price rows into the Cilantro dataset.
# Everything currently at row 322 onward shifts down by two rows, and the
# two newly opened rows (322 and 323) are populated with the new week's
# "caja 36 atados" / "docena de atados" records for 2021-11-04 (serial 44504).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 322.. down by two to make room for the new records.
$ws.Rows.Item(322).Insert()
$ws.Rows.Item(322).Insert()

# New row 322: $/caja 36 atados record
$ws.Cells.Item(322, 1).Value = 9
$ws.Cells.Item(322, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(322, 3).Value = "Metropolitana"
$ws.Cells.Item(322, 4).Value = 44504
$ws.Cells.Item(322, 5).Value = 13
$ws.Cells.Item(322, 6).Value = 100112040
$ws.Cells.Item(322, 7).Value = "Cilantro"
$ws.Cells.Item(322, 8).Value = "Sin especificar"
$ws.Cells.Item(322, 9).Value = "Primera"
$ws.Cells.Item(322, 10).Value = 43
$ws.Cells.Item(322, 11).Value = 4000
$ws.Cells.Item(322, 12).Value = 4500
$ws.Cells.Item(322, 13).Value = 4256
$ws.Cells.Item(322, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(322, 15).Value = "Región Metropolitana"
$ws.Cells.Item(322, 16).Value = 118
$ws.Cells.Item(322, 17).Value = 36
$ws.Cells.Item(322, 18).Value = "Hortaliza"

# New row 323: $/docena de atados record
$ws.Cells.Item(323, 1).Value = 9
$ws.Cells.Item(323, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(323, 3).Value = "Metropolitana"
$ws.Cells.Item(323, 4).Value = 44504
$ws.Cells.Item(323, 5).Value = 13
$ws.Cells.Item(323, 6).Value = 100112040
$ws.Cells.Item(323, 7).Value = "Cilantro"
$ws.Cells.Item(323, 8).Value = "Sin especificar"
$ws.Cells.Item(323, 9).Value = "Primera"
$ws.Cells.Item(323, 10).Value = 160
$ws.Cells.Item(323, 11).Value = 8000
$ws.Cells.Item(323, 12).Value = 10000
$ws.Cells.Item(323, 13).Value = 9000
$ws.Cells.Item(323, 14).Value = "`$/docena de atados"
$ws.Cells.Item(323, 15).Value = "Región Metropolitana"
$ws.Cells.Item(323, 16).Value = 3000
$ws.Cells.Item(323, 17).Value = 3
$ws.Cells.Item(323, 18).Value = "Hortaliza"
